$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 272.4
$ws.Range("I11").Value = 272.4
$ws.Range("K11").Value = 272.4
$ws.Range("M11").Value = -132.4
$ws.Range("H48").Value = 10648.177
$ws.Range("J48").Value = 16002.714
$ws.Range("L48").Value = 48008.142
$ws.Range("N48").Value = -48592.142
$ws.Range("H56").Value = 10648.177
$ws.Range("J56").Value = 16002.714
$ws.Range("L56").Value = 48008.142
$ws.Range("N56").Value = -49076.142
$ws.Range("H62").Value = 6197
$ws.Range("I62").Value = 4080.818
$ws.Range("J62").Value = 9106.75
$ws.Range("K62").Value = 4080.818
$ws.Range("L62").Value = 9106.75
$ws.Range("M62").Value = -3456.818
$ws.Range("N62").Value = -10354.75
$ws.Range("H65").Value = 6197
$ws.Range("I65").Value = 4080.818
$ws.Range("J65").Value = 9106.75
$ws.Range("K65").Value = 20404.09
$ws.Range("L65").Value = 45533.75
$ws.Range("M65").Value = -17284.09
$ws.Range("N65").Value = -51773.75
$ws.Range("H138").Value = 2232.5757
$ws.Range("I138").Value = 809.8421
$ws.Range("J138").Value = 3118.869
$ws.Range("K138").Value = 2429.5263
$ws.Range("L138").Value = 9356.607
$ws.Range("M138").Value = 2710.4737
$ws.Range("N138").Value = -19636.607

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4142.035
$ws.Range("I32").Value = 2732.827
$ws.Range("K32").Value = 2732.827
$ws.Range("M32").Value = -2445.827
$ws.Range("H45").Value = 4050.7144
$ws.Range("I45").Value = 3771
$ws.Range("J45").Value = 4750
$ws.Range("K45").Value = 3771
$ws.Range("L45").Value = 4750
$ws.Range("M45").Value = -3394
$ws.Range("N45").Value = -5504
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("L70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("L73").Value = 0
$ws.Range("H110").Value = 3246.9678
$ws.Range("I110").Value = 2597.889
$ws.Range("K110").Value = 2597.889
$ws.Range("M110").Value = -552.8890000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2853.9644
$ws.Range("I86").Value = 2309.2942
$ws.Range("J86").Value = 3695.7273
$ws.Range("K86").Value = 2309.2942
$ws.Range("L86").Value = 3695.7273
$ws.Range("M86").Value = -1186.2942
$ws.Range("N86").Value = -5941.7273
$ws.Range("H89").Value = 2853.9644
$ws.Range("I89").Value = 2309.2942
$ws.Range("J89").Value = 3695.7273
$ws.Range("K89").Value = 11546.471
$ws.Range("L89").Value = 18478.6365
$ws.Range("M89").Value = -5930.471
$ws.Range("N89").Value = -29710.6365
$ws.Range("H94").Value = 3019
$ws.Range("I94").Value = 2540.3076
$ws.Range("K94").Value = 2540.3076
$ws.Range("M94").Value = -2089.3076
$ws.Range("H107").Value = 1813
$ws.Range("I107").Value = 1813
$ws.Range("K107").Value = 1813
$ws.Range("M107").Value = 107
$ws.Range("H134").Value = 3499.3845
$ws.Range("I134").Value = 3707.75
$ws.Range("K134").Value = 11123.25
$ws.Range("M134").Value = -8588.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 65865.3
$ws.Range("I31").Value = 2558.8
$ws.Range("J31").Value = 92243
$ws.Range("K31").Value = 2558.8
$ws.Range("L31").Value = 92243
$ws.Range("M31").Value = -2263.8
$ws.Range("N31").Value = -92833
$ws.Range("H34").Value = 65865.3
$ws.Range("I34").Value = 2558.8
$ws.Range("J34").Value = 92243
$ws.Range("K34").Value = 2558.8
$ws.Range("L34").Value = 92243
$ws.Range("M34").Value = -2356.8
$ws.Range("N34").Value = -92647
$ws.Range("H58").Value = 6811.353
$ws.Range("I58").Value = 1981
$ws.Range("J58").Value = 15667
$ws.Range("K58").Value = 1981
$ws.Range("L58").Value = 15667
$ws.Range("M58").Value = -1778
$ws.Range("N58").Value = -16073
$ws.Range("H105").Value = 3707.889
$ws.Range("I105").Value = 1393.3334
$ws.Range("J105").Value = 8337
$ws.Range("K105").Value = 1393.3334
$ws.Range("L105").Value = 8337
$ws.Range("M105").Value = 353.6666
$ws.Range("N105").Value = -11831
$ws.Range("H107").Value = 2254.8125
$ws.Range("I107").Value = 1145.8182
$ws.Range("K107").Value = 1145.8182
$ws.Range("M107").Value = 774.1818000000001
$ws.Range("H122").Value = 5563.5293
$ws.Range("I122").Value = 1962
$ws.Range("J122").Value = 12166.333
$ws.Range("K122").Value = 5886
$ws.Range("L122").Value = 36498.999
$ws.Range("M122").Value = -3436
$ws.Range("N122").Value = -41398.999
$ws.Range("H136").Value = 6811.353
$ws.Range("I136").Value = 1981
$ws.Range("J136").Value = 15667
$ws.Range("K136").Value = 5943
$ws.Range("L136").Value = 47001
$ws.Range("M136").Value = -3393
$ws.Range("N136").Value = -52101

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2391.5
$ws.Range("I11").Value = 3274.5
$ws.Range("J11").Value = 1950
$ws.Range("K11").Value = 9823.5
$ws.Range("L11").Value = 5850
$ws.Range("M11").Value = -9683.5
$ws.Range("N11").Value = -6130
$ws.Range("H70").Value = 5663.3335
$ws.Range("I70").Value = 995
$ws.Range("K70").Value = 2985
$ws.Range("M70").Value = -2670
$ws.Range("H73").Value = 5663.3335
$ws.Range("I73").Value = 995
$ws.Range("K73").Value = 2985
$ws.Range("M73").Value = -1893
$ws.Range("H92").Value = 3960
$ws.Range("J92").Value = 4657.143
$ws.Range("L92").Value = 13971.429
$ws.Range("N92").Value = -16467.429
$ws.Range("H131").Value = 13890745
$ws.Range("J131").Value = 8840324
$ws.Range("L131").Value = 26520972
$ws.Range("N131").Value = -26531052
$ws.Range("H132").Value = 4546
$ws.Range("I132").Value = 4166.5
$ws.Range("J132").Value = 4753
$ws.Range("K132").Value = 37498.5
$ws.Range("L132").Value = 42777
$ws.Range("M132").Value = -34968.5
$ws.Range("N132").Value = -47837
$ws.Range("H137").Value = 72028.14
$ws.Range("I137").Value = 672.6667
$ws.Range("J137").Value = 500161
$ws.Range("K137").Value = 2018.0001
$ws.Range("L137").Value = 1500483
$ws.Range("M137").Value = 3081.9999
$ws.Range("N137").Value = -1510683

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 281583.34
$ws.Range("I80").Value = 558153.5600000001
$ws.Range("J80").Value = 5013.1113
$ws.Range("K80").Value = 558153.5600000001
$ws.Range("L80").Value = 5013.1113
$ws.Range("M80").Value = -557155.5600000001
$ws.Range("N80").Value = -7009.1113
$ws.Range("H83").Value = 281583.34
$ws.Range("I83").Value = 558153.5600000001
$ws.Range("J83").Value = 5013.1113
$ws.Range("K83").Value = 2790767.8
$ws.Range("L83").Value = 25065.5565
$ws.Range("M83").Value = -2785775.8
$ws.Range("N83").Value = -35049.5565
$ws.Range("H102").Value = 2798.4
$ws.Range("I102").Value = 2247.95
$ws.Range("K102").Value = 2247.95
$ws.Range("M102").Value = -625.9499999999998

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1340.7
$ws.Range("I16").Value = 812.7059
$ws.Range("K16").Value = 812.7059
$ws.Range("M16").Value = -642.7059
$ws.Range("H22").Value = 3488.65
$ws.Range("I22").Value = 2222.25
$ws.Range("K22").Value = 2222.25
$ws.Range("M22").Value = -1927.25
$ws.Range("H27").Value = 3488.65
$ws.Range("I27").Value = 2222.25
$ws.Range("K27").Value = 2222.25
$ws.Range("M27").Value = -2115.25
$ws.Range("H40").Value = 2004.25
$ws.Range("I40").Value = 1999.5
$ws.Range("J40").Value = 2009
$ws.Range("K40").Value = 1999.5
$ws.Range("L40").Value = 2009
$ws.Range("M40").Value = -1863.5
$ws.Range("N40").Value = -2281
$ws.Range("H61").Value = 4142.364
$ws.Range("I61").Value = 3479.88
$ws.Range("K61").Value = 3479.88
$ws.Range("M61").Value = -3277.88
$ws.Range("H82").Value = 11168.6
$ws.Range("I82").Value = 1412.3334
$ws.Range("J82").Value = 25803
$ws.Range("K82").Value = 1412.3334
$ws.Range("L82").Value = 25803
$ws.Range("M82").Value = -1051.3334
$ws.Range("N82").Value = -26525
$ws.Range("H85").Value = 11168.6
$ws.Range("I85").Value = 1412.3334
$ws.Range("J85").Value = 25803
$ws.Range("K85").Value = 1412.3334
$ws.Range("L85").Value = 25803
$ws.Range("M85").Value = -164.3334
$ws.Range("N85").Value = -28299
$ws.Range("H113").Value = 4142.364
$ws.Range("I113").Value = 3479.88
$ws.Range("K113").Value = 3479.88
$ws.Range("M113").Value = -1309.88
$ws.Range("H122").Value = 676570.3
$ws.Range("I122").Value = 1338666.4
$ws.Range("K122").Value = 4015999.2
$ws.Range("M122").Value = -4013549.2
$ws.Range("H132").Value = 3744.9302
$ws.Range("J132").Value = 4702.8
$ws.Range("L132").Value = 14108.4
$ws.Range("N132").Value = -19168.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9859.799999999999
$ws.Range("I62").Value = 2997.5
$ws.Range("J62").Value = 11575.375
$ws.Range("K62").Value = 2997.5
$ws.Range("L62").Value = 11575.375
$ws.Range("M62").Value = -2373.5
$ws.Range("N62").Value = -12823.375
$ws.Range("H65").Value = 9859.799999999999
$ws.Range("I65").Value = 2997.5
$ws.Range("J65").Value = 11575.375
$ws.Range("K65").Value = 14987.5
$ws.Range("L65").Value = 57876.875
$ws.Range("M65").Value = -11867.5
$ws.Range("N65").Value = -64116.875
$ws.Range("H126").Value = 1388.3334
$ws.Range("I126").Value = 1434.3334
$ws.Range("J126").Value = 1273.3334
$ws.Range("K126").Value = 4303.0002
$ws.Range("L126").Value = 3820.0002
$ws.Range("M126").Value = -1833.0002
$ws.Range("N126").Value = -8760.0002
$ws.Range("H132").Value = 1610.875
$ws.Range("I132").Value = 821.9524
$ws.Range("K132").Value = 2465.8572
$ws.Range("M132").Value = 64.14280000000008

Write-Host "Applied all Lamia_Profits updates"